$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: add the two new header cells (T4, U4) ---
$ws1.Range("T4").Value = "30+"
$ws1.Range("U4").Value = "11x13"

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Column widths for Sheet2
$ws2.Columns.Item(1).ColumnWidth = 13.7265625
$ws2.Columns.Item(2).ColumnWidth = 52.36328125
$ws2.Columns.Item(3).ColumnWidth = 31.90625
$ws2.Columns.Item(4).ColumnWidth = 136.6328125

# Header row
$ws2.Range("A1").Value = "Applied"
$ws2.Range("B1").Value = "Company"
$ws2.Range("C1").Value = "Position"

# Data row
$ws2.Range("A2").Value = (Get-Date -Year 2018 -Month 10 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws2.Range("A2").NumberFormat = "m/d/yyyy"
$ws2.Range("B2").Value = "Precruitment"
$ws2.Range("C2").Value = "Senior Test Analyst"

# --- Selections ---
$ws1.Range("U5").Select()
$ws2.Range("B8").Select()
$ws2.Activate()
